# "assign Remote User condition modified"
#
# Summary of the change (reconstructed from the OOXML diff):
#  - Sheet "Details": a new column is inserted before "SelectAllUserCondition"
#    (old column X), carrying the header "AssignRemoteUser". Everything from
#    old column X onward (SelectAllUserCondition, UsernamesAndRoles,
#    SetNumberOfVirtualMachines, VMName) shifts one column to the right.
#  - Row 2 (the first data row) gets a new "AssignRemoteUser" value ("no")
#    plus a handful of other field edits describing a different test
#    resource (a Windows Server 2022 VM), while rows 3-14 simply shift and
#    leave the brand-new AssignRemoteUser cell blank.
#  - Minor view-state changes: selection on both "ResourceName" and
#    "Details" sheets, plus a couple of column-width tweaks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ResourceName": only the remembered selection changes.
# ---------------------------------------------------------------------
$wsResource = $wb.Worksheets.Item("ResourceName")
$wsResource.Activate()
$wsResource.Range("C2").Select()

# ---------------------------------------------------------------------
# Sheet "Details": the real content edit.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Details")
$ws.Activate()

# Insert a new column at X (pushes SelectAllUserCondition/UsernamesAndRoles/
# SetNumberOfVirtualMachines/VMName one column to the right, and grows the
# used range from A1:AA14 to A1:AB14 automatically).
$ws.Range("X1").EntireColumn.Insert()

# Row 2 edits - create the new shared strings in the same order they first
# appear left-to-right in the edited row so the workbook's string table
# matches the target layout.
$ws.Range("E2").Value = "Windows Server 2022 (Standard Edition without License)"
$ws.Range("W2").Value = "10.150.35.0 - akashTest"

# New column header (new shared string created last).
$ws.Range("X1").Value = "AssignRemoteUser"

# Remaining row 2 field edits.
$ws.Range("C2").Value = "EC_Windows_VM"
$ws.Range("G2").Value = "none"
$ws.Range("L2").Value = "4"
$ws.Range("M2").Value = "4"
$ws.Range("P2").Value = "yes"
$ws.Range("S2").Value = "Weekly"
$ws.Range("T2").Value = "6 weeks"
$ws.Range("U2").Value = "yes"

# New "AssignRemoteUser" value for row 2 only; rows 3-14 are left blank in
# this column.
$ws.Range("X2").Value = "no"

# Column width tweaks (column O widened, new column X given an explicit
# width).
$ws.Columns.Item(15).ColumnWidth = 27.333333333333332
$ws.Columns.Item(24).ColumnWidth = 16.833333333333332

# Restore the remembered selection on this sheet, and make sure it stays the
# active tab (it was, and still is, the active sheet in the workbook).
$ws.Range("Z13").Select()
